$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "last refreshed" timestamp string in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 16:05"

# --- 2. Refresh per-country COVID numbers (columns B..H) ---
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1529291
$ws.Cells.Item(4, 3).Value = 1627
$ws.Cells.Item(4, 4).Value = 346394
$ws.Cells.Item(4, 5).Value = 1091892
$ws.Cells.Item(4, 7).Value = 27
$ws.Cells.Item(4, 8).Value = 91005

# Row 7: Brasil
$ws.Cells.Item(7, 2).Value = 244052
$ws.Cells.Item(7, 3).Value = 2972
$ws.Cells.Item(7, 5).Value = 133729
$ws.Cells.Item(7, 7).Value = 83
$ws.Cells.Item(7, 8).Value = 16201

# Row 52: Noruega
$ws.Cells.Item(52, 2).Value = 8257
$ws.Cells.Item(52, 3).Value = 8
$ws.Cells.Item(52, 5).Value = 7992

# Row 53: Argentina
$ws.Cells.Item(53, 4).Value = 2625
$ws.Cells.Item(53, 5).Value = 5069
$ws.Cells.Item(53, 7).Value = 1
$ws.Cells.Item(53, 8).Value = 374

# Row 71: Azerbaiyan
$ws.Cells.Item(71, 2).Value = 3387
$ws.Cells.Item(71, 3).Value = 113
$ws.Cells.Item(71, 4).Value = 2055
$ws.Cells.Item(71, 5).Value = 1292
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 40

# Row 80: Bosnia y Herzegovina
$ws.Cells.Item(80, 2).Value = 2304
$ws.Cells.Item(80, 3).Value = 14
$ws.Cells.Item(80, 4).Value = 1464
$ws.Cells.Item(80, 5).Value = 707

# Row 85: Cuba
$ws.Cells.Item(85, 2).Value = 1881
$ws.Cells.Item(85, 3).Value = 9
$ws.Cells.Item(85, 4).Value = 1505
$ws.Cells.Item(85, 5).Value = 297

# --- 3. Reorder "Santa Lucia" ahead of "Belice" / "Nueva Caledonia" ---
# (rows 195 = Belice, 196 = Nueva Caledonia, 197 = Santa Lucia today;
#  Santa Lucia should move up to sit right after Fiyi at row 195,
#  pushing Belice and Nueva Caledonia down by one row each)
$row195 = @()
$row196 = @()
$row197 = @()
for ($c = 1; $c -le 8; $c++) {
    $row195 += ,$ws.Cells.Item(195, $c).Value()
    $row196 += ,$ws.Cells.Item(196, $c).Value()
    $row197 += ,$ws.Cells.Item(197, $c).Value()
}

for ($c = 1; $c -le 8; $c++) {
    $ws.Cells.Item(195, $c).Value = $row197[$c - 1]
    $ws.Cells.Item(196, $c).Value = $row195[$c - 1]
    $ws.Cells.Item(197, $c).Value = $row196[$c - 1]
}
